$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '34.479.25'
$ws.Range('E2').Value = '  +0.13%  '
Set-TextValue 'D3' '1.807.01'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '224.92'
$ws.Range('E5').Value = '  -1.18%  '
Set-TextValue 'D6' '0.600'
$ws.Range('E6').Value = '  +4.66%  '
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue 'D8' '38.59'
$ws.Range('E8').Value = '  +6.58%  '
$ws.Range('E9').Value = '  -3.68%  '
Set-TextValue 'D10' '0.0669'
$ws.Range('E10').Value = '  -3.12%  '
Set-TextValue 'D11' '0.0981'
$ws.Range('E11').Value = '  +1.86%  '
Set-TextValue 'D12' '2.068.74'
$ws.Range('E12').Value = '  +0.38%  '
Set-TextValue 'D13' '11.10'
$ws.Range('E13').Value = '  -4.16%  '
Set-TextValue 'D14' '1.797.84'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -2.22%  '
Set-TextValue 'D16' '34.470.56'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('E17').Value = '  -2.55%  '
Set-TextValue 'D18' '68.13'
$ws.Range('E18').Value = '  -1.25%  '
Set-TextValue 'D19' '240.84'
$ws.Range('E19').Value = '  -1.63%  '
$ws.Range('E20').Value = '  -2.74%  '
Set-TextValue 'D21' '11.17'
$ws.Range('E21').Value = '  -3.61%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('E24').Value = '  +0.82%  '
Set-TextValue 'D25' '170.66'
$ws.Range('E25').Value = '  -1.05%  '
Set-TextValue 'D26' '7.68'
$ws.Range('E26').Value = '  -3.61%  '
Set-TextValue 'D27' '17.46'
$ws.Range('E27').Value = '  +3.76%  '
Set-TextValue 'D28' '0.121'
$ws.Range('E28').Value = '  +3.24%  '
$ws.Range('E29').Value = '  +0.05%  '
Set-TextValue 'D30' '1.23'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('E31').Value = '  -1.74%  '
$ws.Range('E32').Value = '  -2.71%  '
$ws.Range('E33').Value = '  -4.24%  '
$ws.Range('E34').Value = '  +0.80%  '
Set-TextValue 'D35' '0.640'
$ws.Range('E35').Value = '  -4.81%  '
Set-TextValue 'D36' '1.305.04'
$ws.Range('E36').Value = '  -6.56%  '
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('E38').Value = '  -1.84%  '
$ws.Range('E39').Value = '  -4.69%  '
Set-TextValue 'D40' '82.69'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D41' '2.44'
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D42' '1.22'
$ws.Range('E42').Value = '  +3.34%  '
$ws.Range('E43').Value = '  -0.56%  '
Set-TextValue 'D44' '0.951'
$ws.Range('E44').Value = '  -0.73%  '
Set-TextValue 'D45' '14.00'
$ws.Range('E45').Value = '  +4.96%  '
Set-TextValue 'D47' '1.969.48'
$ws.Range('E47').Value = '  +0.32%  '
Set-TextValue 'D48' '5.80'
$ws.Range('E48').Value = '  -3.83%  '
$ws.Range('E49').Value = '  -0.03%  '
Set-TextValue 'D50' '102.86'
$ws.Range('E50').Value = '  -1.44%  '
$ws.Range('E51').Value = '  -1.81%  '
